# Auto-generated edit script: updates crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.912.02"
$ws.Range("E2").Value = "  -3.60%  "
$ws.Range("D3").Value = "2.920.37"
$ws.Range("E3").Value = "  -4.06%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'586.51"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'145.39"
$ws.Range("E6").Value = "  -6.14%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.505"
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("D9").Value = "2.919.58"
$ws.Range("E9").Value = "  -4.11%  "
$ws.Range("D10").Value = "'6.83"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").Value = "'0.144"
$ws.Range("E11").Value = "  -4.97%  "
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("E13").Value = "  -4.08%  "
$ws.Range("D14").Value = "'33.64"
$ws.Range("E14").Value = "  -6.37%  "
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "3.404.53"
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("D17").Value = "60.901.40"
$ws.Range("E17").Value = "  -3.51%  "
$ws.Range("D18").Value = "'6.77"
$ws.Range("E18").Value = "  -4.82%  "
$ws.Range("D19").Value = "2.922.05"
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("D20").Value = "'429.36"
$ws.Range("E20").Value = "  -5.93%  "
$ws.Range("D21").Value = "'13.62"
$ws.Range("E21").Value = "  -5.08%  "
$ws.Range("D22").Value = "'0.683"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "'7.13"
$ws.Range("E23").Value = "  -5.55%  "
$ws.Range("D24").Value = "'80.68"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("D25").Value = "'2.25"
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("D26").Value = "'10.75"
$ws.Range("E26").Value = "  -4.48%  "
$ws.Range("D27").Value = "'11.99"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "'7.20"
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("D35").Value = "0.0₃0871"
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("D36").Value = "'1.02"
$ws.Range("E36").Value = "  -2.89%  "
$ws.Range("D37").Value = "'5.66"
$ws.Range("E37").Value = "  -5.06%  "
$ws.Range("D38").Value = "'3.02"
$ws.Range("E38").Value = "  -5.75%  "
$ws.Range("E39").Value = "  -3.96%  "
$ws.Range("D40").Value = "'49.66"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").Value = "'2.00"
$ws.Range("E41").Value = "  -5.81%  "
$ws.Range("D42").Value = "'8.66"
$ws.Range("E42").Value = "  -5.41%  "
$ws.Range("D43").Value = "'0.299"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").Value = "'40.97"
$ws.Range("E44").Value = "  -6.41%  "
$ws.Range("D45").Value = "'0.0353"
$ws.Range("E45").Value = "  -2.88%  "
$ws.Range("D46").Value = "'379.48"
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("D47").Value = "2.692.84"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").Value = "'132.74"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D50").Value = "'24.51"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -2.36%  "

# Reset style on cells that needed a text-forcing apostrophe prefix,
# so no stray quote-prefix style attribute lingers on the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
